# Trading update: 2026-02-17 20:27:45
# Appends the newest OPEN MarketMaking trade (Trade # 52) to both the
# "All Trades" log and the per-strategy "MarketMaking" sheet.

$wb = $excel.ActiveWorkbook

$tradeNumber   = 52
$tradeDate     = "2026-02-17"
$tradeTime     = "20:27:06"
$strategy      = "MarketMaking"
$side          = "DOWN"
$entryPrice    = 0.47
$status        = "OPEN"
$pnlPct        = 0
$pnlDollar     = 0
$capitalAfter  = 100
$entrySlippage = 0
$exitSlippage  = 0
$confidence    = 0.6
$entryReason   = "Normal spread capture: 19600 bps"
$duration      = 0

function Add-TradeRow($ws, $row) {
    $ws.Cells.Item($row, 1).Value = $tradeNumber
    $ws.Cells.Item($row, 2).Value = $tradeDate
    $ws.Cells.Item($row, 3).Value = $tradeTime
    $ws.Cells.Item($row, 4).Value = $strategy
    $ws.Cells.Item($row, 5).Value = $side
    $ws.Cells.Item($row, 6).Value = $entryPrice
    $ws.Cells.Item($row, 7).Value = ""
    $ws.Cells.Item($row, 8).Value = $status
    $ws.Cells.Item($row, 9).Value = $pnlPct
    $ws.Cells.Item($row, 10).Value = $pnlDollar
    $ws.Cells.Item($row, 11).Value = $capitalAfter
    $ws.Cells.Item($row, 12).Value = $entrySlippage
    $ws.Cells.Item($row, 13).Value = $exitSlippage
    $ws.Cells.Item($row, 14).Value = $confidence
    $ws.Cells.Item($row, 15).Value = $entryReason
    $ws.Cells.Item($row, 16).Value = ""
    $ws.Cells.Item($row, 17).Value = $duration
}

$allTrades = $wb.Worksheets.Item("All Trades")
Add-TradeRow $allTrades 53

$marketMaking = $wb.Worksheets.Item("MarketMaking")
Add-TradeRow $marketMaking 20
